$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.392.72"
$ws.Range("E2").Value = "  +4.20%  "
$ws.Range("D3").Value = "2.355.31"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'546.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("D6").Value = "'132.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").Value = "2.350.52"
$ws.Range("E9").Value = "  +2.74%  "
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("D11").Value = "'5.52"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "'23.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").Value = "2.768.78"
$ws.Range("E15").Value = "  +2.62%  "
$ws.Range("D16").Value = "60.401.31"
$ws.Range("E16").Value = "  +4.31%  "
$ws.Range("D17").Value = "'0.0000134"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.93%  "
$ws.Range("D18").Value = "2.352.98"
$ws.Range("E18").Value = "  +2.49%  "
$ws.Range("D19").Value = "'10.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("D20").Value = "'4.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").Value = "'6.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.44%  "
$ws.Range("D22").Value = "'315.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("E25").Value = "  +2.50%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("E28").Value = "  +5.28%  "
$ws.Range("E29").Value = "  +3.63%  "
$ws.Range("D30").Value = "'171.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("E31").Value = "  +10.29%  "
$ws.Range("D32").Value = "0.0₃0732"
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("D33").Value = "'5.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.20%  "
$ws.Range("D34").Value = "'1.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +14.84%  "
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("D36").Value = "'18.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.86%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  +6.98%  "
$ws.Range("D40").Value = "'317.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.97%  "
$ws.Range("D41").Value = "'38.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("D42").Value = "'1.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.49%  "
$ws.Range("D43").Value = "'142.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("E44").Value = "  +1.53%  "
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("D46").Value = "'19.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.95%  "
$ws.Range("E47").Value = "  +1.03%  "
$ws.Range("D48").Value = "'0.562"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.50%  "
$ws.Range("E49").Value = "  +2.11%  "
$ws.Range("D50").Value = "'11.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("D51").Value = "0.0₆0210"
$ws.Range("E51").Value = "  +10.91%  "
